$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70; existing rows 70-147 shift down to 71-148
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record
$ws.Cells.Item(70, 1).Value = 1
$ws.Cells.Item(70, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(70, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(70, 4).Value = 44904
$ws.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70, 5).Value = 15
$ws.Cells.Item(70, 6).Value = 100112042
$ws.Cells.Item(70, 7).Value = "Locoto"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 100
$ws.Cells.Item(70, 11).Value = 18000
$ws.Cells.Item(70, 12).Value = 19000
$ws.Cells.Item(70, 13).Value = 18500
$ws.Cells.Item(70, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 925
$ws.Cells.Item(70, 17).Value = 20
$ws.Cells.Item(70, 18).Value = "Hortaliza"
